$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 124
$ws.Range("I2").Value = 316
$ws.Range("J2").Value = 1296
$ws.Range("L2").Value = 345
$ws.Range("N2").Value = 231
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 18
$ws.Range("S2").Value = 137
$ws.Range("T2").Value = 228
$ws.Range("U2").Value = 20
$ws.Range("V2").Value = 1893
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 2005
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 26
$ws.Range("AA2").Value = 15
